$d = $word.ActiveDocument

function Replace-Exact([string]$old, [string]$new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 1. {{LEITBEHOERDE_STADT}}, {{HEUTE}}  (merge split runs, was split across
#    "{{LEITBEHOERDE_" / "STADT" / "}}, {{" / "HEUTE" / "}}")
Replace-Exact "{{LEITBEHOERDE_STADT}}, {{HEUTE}}" "{{LEITBEHOERDE_STADT}}, {{HEUTE}}"

# 2. {{LEITBEHOERDE_ADRESSE_1  (merge "{{LEITBEHOERDE_ADRESS" / "E" / "_1",
#    keep the bookmark-interrupted trailing "}}" untouched)
Replace-Exact "{{LEITBEHOERDE_ADRESSE_1" "{{LEITBEHOERDE_ADRESSE_1"

# 3. {{LEITBEHOERDE_ADRESSE_2}}  (merge "{{LEITBEHOERDE_ADRESS" / "E" / "_2}}")
Replace-Exact "{{LEITBEHOERDE_ADRESSE_2}}" "{{LEITBEHOERDE_ADRESSE_2}}"

# 4. Telefon {{LEITBEHOERDE_TELEFON}}
Replace-Exact "Telefon {{LEITBEHOERDE_TELEFON}}" "Telefon {{LEITBEHOERDE_TELEFON}}"

# 5. {{ADRESSE}}
Replace-Exact "{{ADRESSE}}" "{{ADRESSE}}"

# 6. {{GESUCHSTELLER_NAME_ADRESSE}}
Replace-Exact "{{GESUCHSTELLER_NAME_ADRESSE}}" "{{GESUCHSTELLER_NAME_ADRESSE}}"

# 7. {{PROJEKTVERFASSER_NAME_ADRESSE}}
Replace-Exact "{{PROJEKTVERFASSER_NAME_ADRESSE}}" "{{PROJEKTVERFASSER_NAME_ADRESSE}}"

# 8. {{FACHSTELLEN_KANTONAL_LISTE}}
Replace-Exact "{{FACHSTELLEN_KANTONAL_LISTE}}" "{{FACHSTELLEN_KANTONAL_LISTE}}"

# 9. Fix placeholder collision: _TYPE}} -> _TYP}}
Replace-Exact "_TYPE}}" "_TYP}}"

# 10. Merge "Bau- " / "und" / " Verkehrsdirektion..." into a single run
Replace-Exact "Dieser Entscheid kann innert 30 Tagen seit der Eröffnung mit Baubeschwerde bei der Bau- und Verkehrsdirektion des Kantons Bern, Reiterstrasse 11, 3011 Bern, angefochten werden (Art. 40 BauG)." "Dieser Entscheid kann innert 30 Tagen seit der Eröffnung mit Baubeschwerde bei der Bau- und Verkehrsdirektion des Kantons Bern, Reiterstrasse 11, 3011 Bern, angefochten werden (Art. 40 BauG)."
